# Append two new daily rows (dates 45985 / 45986) to each of the six
# worksheets, matching the "date" / "remn_amt" column layout already in
# the sheet. New date cells reuse the same date number-format/style as
# the existing date column; the amount column keeps the default (General)
# style, just like every other row.

$wb = $excel.ActiveWorkbook

$dateNumberFormat = "YYYY-MM-DD HH:MM:SS"

# sheet index -> starting row for the new data, and the two (date, amount) pairs to add
$newRows = @{
    1 = @{ StartRow = 103; Values = @(@(45985, 662265), @(45986, 612245)) }
    2 = @{ StartRow = 103; Values = @(@(45985, 2019249), @(45986, 2064804)) }
    3 = @{ StartRow = 103; Values = @(@(45985, 337342), @(45986, 349928)) }
    4 = @{ StartRow = 103; Values = @(@(45985, 151232), @(45986, 145925)) }
    5 = @{ StartRow = 68;  Values = @(@(45985, 24874),  @(45986, 35487)) }
    6 = @{ StartRow = 103; Values = @(@(45985, 48808),  @(45986, 48608)) }
}

foreach ($sheetIndex in 1..6) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $info = $newRows[$sheetIndex]
    $row = $info.StartRow

    foreach ($pair in $info.Values) {
        $dateVal = $pair[0]
        $amountVal = $pair[1]

        $dateCell = $ws.Cells.Item($row, 1)
        $dateCell.Value = $dateVal
        $dateCell.NumberFormat = $dateNumberFormat

        $amountCell = $ws.Cells.Item($row, 2)
        $amountCell.Value = $amountVal

        $row = $row + 1
    }
}
